$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Column C (Yes -> No) values in existing rows 5-7 ---
$ws.Range("C5").Value = "No"
$ws.Range("C6").Value = "No"

# Row 7's C cell did not previously carry the shared "left/top" style (s=1);
# copy that formatting down from C6 before setting its value.
$ws.Range("C6").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = "No"

# BX5 previously carried an extra fill style; the new layout no longer needs it.
$ws.Range("BX5").ClearFormats()

# --- Add new rows 8-10 for the LAF test cases ---
$ws.Range("A8").Value = "TC_07_Validate_LAF_Title"
$ws.Range("B8").Value = "Membership"
$ws.Range("C3").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C8").Value = "Yes"
$ws.Range("F8").Value = "CheckLAF_Title"
$ws.Range("BX8").Value = "LA Fitness | Gym and Fitness Club | Join Today"

$ws.Range("A9").Value = "TC_08_Validate_JoinNow_Button"
$ws.Range("B9").Value = "Membership"
$ws.Range("C3").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C9").Value = "Yes"
$ws.Range("F9").Value = "Check_LAF_JoinNowHeaderButton"

$ws.Range("A10").Value = "TC_09_Validate_JoinNow_ClickButton"
$ws.Range("B10").Value = "Membership"
$ws.Range("C3").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("C10").Value = "Yes"
$ws.Range("F10").Value = "Click_JoinNowHeaderButton"

# --- Move the active selection to match the saved view state ---
$ws.Range("F17").Select() | Out-Null
